$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their text formatting so values
# like "246.50" or "37.167.29" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "37.167.29"
$ws.Range("E2").Value = "  +1.80%  "
$ws.Range("D3").Value = "1.999.99"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "246.50"
$ws.Range("E5").Value = "  +0.92%  "
$ws.Range("D6").Value = "0.629"
$ws.Range("E6").Value = "  +2.53%  "
$ws.Range("D7").Value = "60.11"
$ws.Range("E7").Value = "  +3.15%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").Value = "0.384"
$ws.Range("E9").Value = "  +1.79%  "
$ws.Range("D10").Value = "0.0802"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "14.89"
$ws.Range("E12").Value = "  +5.45%  "
$ws.Range("D13").Value = "22.60"
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("D14").Value = "2.293.57"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").Value = "0.845"
$ws.Range("E15").Value = "  +0.46%  "
$ws.Range("E16").Value = "  +3.14%  "
$ws.Range("D17").Value = "1.999.56"
$ws.Range("E17").Value = "  +2.47%  "
$ws.Range("D18").Value = "37.079.22"
$ws.Range("E18").Value = "  +1.58%  "
$ws.Range("D19").Value = "70.27"
$ws.Range("E19").Value = "  +0.80%  "
$ws.Range("D20").Value = "0.0₃0863"
$ws.Range("E20").Value = "  +1.74%  "
$ws.Range("D21").Value = "5.17"
$ws.Range("E21").Value = "  +2.63%  "
$ws.Range("D22").Value = "230.59"
$ws.Range("E22").Value = "  +0.51%  "
$ws.Range("E23").Value = "  +0.19%  "
$ws.Range("D24").Value = "2.47"
$ws.Range("E24").Value = "  +1.23%  "
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "9.41"
$ws.Range("E26").Value = "  +2.91%  "
$ws.Range("D27").Value = "0.142"
$ws.Range("E27").Value = "  +4.35%  "
$ws.Range("D28").Value = "163.55"
$ws.Range("E28").Value = "  +1.96%  "
$ws.Range("D29").Value = "19.62"
$ws.Range("E29").Value = "  +1.47%  "
$ws.Range("D30").Value = "1.35"
$ws.Range("E30").Value = "  +14.87%  "
$ws.Range("E31").Value = "  +1.20%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  +7.61%  "
$ws.Range("D34").Value = "4.53"
$ws.Range("E34").Value = "  +3.16%  "
$ws.Range("D35").Value = "2.38"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("D37").Value = "1.80"
$ws.Range("E37").Value = "  +2.21%  "
$ws.Range("D38").Value = "3.26"
$ws.Range("E38").Value = "  -4.56%  "
$ws.Range("D39").Value = "5.42"
$ws.Range("E39").Value = "  +3.67%  "
$ws.Range("E40").Value = "  +0.25%  "
$ws.Range("E41").Value = "  +0.46%  "
$ws.Range("D42").Value = "0.0215"
$ws.Range("E42").Value = "  +2.27%  "
$ws.Range("E43").Value = "  +1.66%  "
$ws.Range("E44").Value = "  +5.04%  "
$ws.Range("D45").Value = "90.94"
$ws.Range("E45").Value = "  +3.56%  "
$ws.Range("D46").Value = "1.374.75"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E47").Value = "  +1.90%  "
$ws.Range("D48").Value = "7.28"
$ws.Range("E48").Value = "  +1.65%  "
$ws.Range("E49").Value = "  +0.51%  "
$ws.Range("E50").Value = "  +14.33%  "
$ws.Range("D51").Value = "46.33"
$ws.Range("E51").Value = "  +5.38%  "
